$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report title text (HS = "Ho so" instead of "Doanh so")
$ws.Range("A1").Value = "BÁO CÁO TỔNG HỢP - TỔNG HỒ SƠ CHI TRẢ"
$ws.Range("Q7").Value = "Tổng Hồ sơ chi trả"

# Update the sheet view: clear frozen/top-left cell, change selection range
$ws.Range("A1:U1").Select()
